$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 63
$ws.Range("B2").Value = 322
$ws.Range("C2").Value = 549
$ws.Range("D2").Value = 179
$ws.Range("E2").Value = 713
